# Add a new worksheet "CustomerDetails" at the end of the workbook (after
# "MultipleData"), populate it with customer/tester data, and style it to
# match the existing "MultipleData" sheet's look (yellow header row with
# borders, bordered data rows, autosized columns, selection on C10).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "CustomerDetails"

# --- Data ---------------------------------------------------------------
$newSheet.Range("A1").Value = "CustomerName"
$newSheet.Range("B1").Value = "CustomerDescription"

$newSheet.Range("A2").Value = "Krishna"
$newSheet.Range("B2").Value = "Tester"

$newSheet.Range("A3").Value = "pakeera"
$newSheet.Range("B3").Value = "Test Lead"

$newSheet.Range("A4").Value = "Vamsi"
$newSheet.Range("B4").Value = "Devloer"

$newSheet.Range("A5").Value = "Shyam"
$newSheet.Range("B5").Value = "ManualTester"

# --- Styling (reuse the existing MultipleData look: header row + bordered
# data rows) by copying the formatting from that sheet so the same
# style/border/fill records are reused instead of new ones being created.
$srcSheet = $wb.Worksheets.Item("MultipleData")

$srcSheet.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

$srcSheet.Range("A2:B2").Copy()
$newSheet.Range("A2:B5").PasteSpecial(-4122)

# --- Column widths --------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 13
$newSheet.Columns.Item(2).ColumnWidth = 17.5

# --- Selection ------------------------------------------------------------
$newSheet.Range("C10").Select()

Write-Output "done"
